$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medium traffic density")

# Filtered out duplicate CPA events: update NMACs ("I" column) values
$ws.Range("I12").Value = 9
$ws.Range("I17").Value = 15
$ws.Range("I22").Value = 19
$ws.Range("I24").Value = 23

# Move active selection to I12 (matches saved cursor position in diff)
$ws.Range("I12").Select()
